$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - these values are numeric-looking text stored as
# inline strings in the source data, so each one is prefixed with a leading
# apostrophe here to force Excel to keep storing it as Text (matching the
# original cell type) rather than silently converting it to a Number cell.
$ws.Range("D2").Value  = "'273.91"
$ws.Range("D3").Value  = "'21.14"
$ws.Range("D4").Value  = "'6.206"
$ws.Range("D5").Value  = "'0.06173"
$ws.Range("D6").Value  = "'3.574"
$ws.Range("D7").Value  = "'1.515"
$ws.Range("D8").Value  = "'6.534"
$ws.Range("D9").Value  = "'0.8230"
$ws.Range("D10").Value = "'0.1648"
$ws.Range("D11").Value = "'0.08213"
$ws.Range("D12").Value = "'0.03435"
$ws.Range("D13").Value = "'0.03145"
$ws.Range("D14").Value = "'0.09131"
$ws.Range("D15").Value = "'3.774"
$ws.Range("D16").Value = "'0.001604"
$ws.Range("D17").Value = "'0.04691"
$ws.Range("D18").Value = "'0.006450"
$ws.Range("D24").Value = "'0.01389"
$ws.Range("D25").Value = "'0.3340"
$ws.Range("D26").Value = "'0.1231"
$ws.Range("D42").Value = "'0.007034"
$ws.Range("D43").Value = "'0.1104"
$ws.Range("D44").Value = "'0.01025"
$ws.Range("D45").Value = "'0.00006563"
$ws.Range("D47").Value = "'0.7231"

# Volume(1h) label (column E) updates: the "Best in 24h" tag moved from row 41
# (CEJI) to row 42 (KickToken).
$ws.Range("E41").Value = "40CEJICEJI"
$ws.Range("E42").Value = "41KickTokenKICKBestin24h"
